# Apply crypto price/volume updates as described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain decimal number (e.g. "3.36") must be forced
# to Text first, otherwise Excel COM auto-converts Range.Value into a numeric
# type (losing the original text formatting / exact digits).
$textForceCells = @(
    "D5",
    "D6",
    "D8",
    "D9",
    "D11",
    "D13",
    "D16",
    "D17",
    "D18",
    "D21",
    "D23",
    "D25",
    "D27",
    "D28",
    "D29",
    "D31",
    "D32",
    "D34",
    "D35",
    "D36",
    "D38",
    "D39",
    "D40",
    "D41",
    "D42",
    "D45",
    "D47",
    "D51",
)
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# --- Plain text / non-numeric-looking values (safe to assign directly) ---
$ws.Range("D2").Value = "73.100.53"
$ws.Range("E2").Value = "  +2.94%  "
$ws.Range("D3").Value = "3.994.93"
$ws.Range("E3").Value = "  +1.24%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("E5").Value = "  +10.95%  "
$ws.Range("E6").Value = "  +11.53%  "
$ws.Range("E7").Value = "  -0.29%  "
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("E9").Value = "  +1.70%  "
$ws.Range("E10").Value = "  +2.12%  "
$ws.Range("E11").Value = "  -0.83%  "
$ws.Range("E13").Value = "  +3.78%  "
$ws.Range("D14").Value = "4.636.93"
$ws.Range("E14").Value = "  +1.43%  "
$ws.Range("D15").Value = "4.001.90"
$ws.Range("E15").Value = "  +1.44%  "
$ws.Range("E16").Value = "  +8.69%  "
$ws.Range("E17").Value = "  +1.76%  "
$ws.Range("E18").Value = "  +0.08%  "
$ws.Range("E19").Value = "  +0.55%  "
$ws.Range("D20").Value = "72.819.01"
$ws.Range("E20").Value = "  +2.80%  "
$ws.Range("E21").Value = "  +4.13%  "
$ws.Range("E22").Value = "  +11.59%  "
$ws.Range("E23").Value = "  -0.83%  "
$ws.Range("E24").Value = "  -4.25%  "
$ws.Range("E25").Value = "  -0.32%  "
$ws.Range("E26").Value = "  +14.99%  "
$ws.Range("E27").Value = "  +0.38%  "
$ws.Range("E28").Value = "  +0.79%  "
$ws.Range("E29").Value = "  -2.76%  "
$ws.Range("E30").Value = "  +0.03%  "
$ws.Range("E31").Value = "  +1.25%  "
$ws.Range("E32").Value = "  +4.45%  "
$ws.Range("E33").Value = "  -0.31%  "
$ws.Range("E34").Value = "  -5.73%  "
$ws.Range("E35").Value = "  -0.63%  "
$ws.Range("E36").Value = "  +7.94%  "
$ws.Range("D37").Value = "0.0₃0905"
$ws.Range("E37").Value = "  +11.11%  "
$ws.Range("E38").Value = "  -0.55%  "
$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("E39").Value = "  -1.11%  "
$ws.Range("B40").Value = "ThetaToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("E40").Value = "  -0.45%  "
$ws.Range("B41").Value = "Dai"
$ws.Range("C41").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("E41").Value = "  -0.02%  "
$ws.Range("B42").Value = "WEMIXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("E42").Value = "  +5.86%  "
$ws.Range("E43").Value = "  +0.44%  "
$ws.Range("E44").Value = "  +2.02%  "
$ws.Range("E45").Value = "  +5.69%  "
$ws.Range("E46").Value = "  +0.60%  "
$ws.Range("E47").Value = "  -1.49%  "
$ws.Range("E48").Value = "  +1.53%  "
$ws.Range("D49").Value = "2.914.69"
$ws.Range("E49").Value = "  +12.60%  "
$ws.Range("E50").Value = "  +1.77%  "
$ws.Range("E51").Value = "  +4.79%  "

# --- Values that needed the Text number-format guard above ---
$ws.Range("D5").Value = "597.81"
$ws.Range("D6").Value = "164.53"
$ws.Range("D8").Value = "0.999"
$ws.Range("D9").Value = "0.752"
$ws.Range("D11").Value = "55.03"
$ws.Range("D13").Value = "11.01"
$ws.Range("D16").Value = "1.26"
$ws.Range("D17").Value = "14.13"
$ws.Range("D18").Value = "20.52"
$ws.Range("D21").Value = "438.50"
$ws.Range("D23").Value = "96.75"
$ws.Range("D25").Value = "14.37"
$ws.Range("D27").Value = "11.39"
$ws.Range("D28").Value = "5.95"
$ws.Range("D29").Value = "10.40"
$ws.Range("D31").Value = "7.95"
$ws.Range("D32").Value = "13.91"
$ws.Range("D34").Value = "48.24"
$ws.Range("D35").Value = "670.89"
$ws.Range("D36").Value = "70.66"
$ws.Range("D38").Value = "0.438"
$ws.Range("D39").Value = "0.146"
$ws.Range("D40").Value = "3.36"
$ws.Range("D41").Value = "0.998"
$ws.Range("D42").Value = "3.35"
$ws.Range("D45").Value = "10.72"
$ws.Range("D47").Value = "2.63"
$ws.Range("D51").Value = "3.42"

# Restore the default (unstyled) cell style now that the text value is locked in,
# so we do not leave a stray NumberFormat behind on these cells.
foreach ($addr in $textForceCells) {
    $ws.Range($addr).Style = "Normal"
}
